$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above row 4, pushing existing rows 4-6 down to 5-7.
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with a full data row (weekly price update).
$ws.Cells.Item(4, 1).Value = 1
$ws.Cells.Item(4, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(4, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(4, 4).Value = 44883
$ws.Cells.Item(4, 5).Value = 15
$ws.Cells.Item(4, 6).Value = 100112030
$ws.Cells.Item(4, 7).Value = "Poroto granado"
$ws.Cells.Item(4, 8).Value = "Sin especificar"
$ws.Cells.Item(4, 9).Value = "Primera"
$ws.Cells.Item(4, 10).Value = 290
$ws.Cells.Item(4, 11).Value = 1400
$ws.Cells.Item(4, 12).Value = 1500
$ws.Cells.Item(4, 13).Value = 1434
$ws.Cells.Item(4, 14).Value = "$/kilo"
$ws.Cells.Item(4, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(4, 16).Value = 1434
$ws.Cells.Item(4, 17).Value = 1
$ws.Cells.Item(4, 18).Value = "Hortaliza"
